# Manual daily update: append 10 new WIP contract rows (362-371) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New contract rows, in column order A..J:
# Date Awarded | Contract # | Qty | Contract Total | NSN | Part Name | Vendor Name | P/N | Preservation Method | Due Date
$newRows = @(
    @('9/19/2019', 'SPE7M1-19-V-057Q', '13', '$2,639.00', '5935015360017', 'DUMMY CONNECTOR,RECEPTACLE', 'Glenair', '650HS058NF25', 'CP', '2020 MAR 17'),
    @('9/12/2019', 'SPE7M5-19-V-116M', '5', '$1,915.00', '5935016512754', 'BACKSHELL,ELECTRICA', 'Glenair', '443AS024NF1404-3S', 'M31', '2020 MAR 10'),
    @('9/12/2019', 'SPE7M1-19-V-046H', '5', '$849.45', '5935014667112', 'SHELL,ELECTRICAL CO', 'Kern', 'BAT86-19W07C', 'CP', '2020 FEB 19'),
    @('9/12/2019', 'SPE7M5-19-V-116S', '6', '$805.20', '5935015754171', 'COVER,ELECTRICAL CONNECTOR', 'Glenair', '660-008NF18F6-12', 'ZZ', '2020 MAR 10'),
    @('9/12/2019', 'SPE7M5-19-V-116T', '5', '$1,643.05', '5935016499431', 'BACKSHELL,ELECTRICA', 'Glenair', '440AS156NF2216-16-4', 'M41', '2020 MAR 10'),
    @('9/12/2019', 'SPE7M5-19-V-115D', '146', '$7,089.76', '5935015812863', 'BACKSHELL,ELECTRICAL CONNECTOR', 'Glenair', '447HS105NF0901K', 'CP', '2020 MAR 10'),
    @('9/12/2019', 'SPE7M5-19-V-116K', '12', '$1,244.28', '5935015754188', 'COVER,ELECTRICAL CO', 'Glenair', '660-008NF22F6-16', 'CP', '2020 MAR 10'),
    @('9/20/2019', 'SPE5E3-19-V-061E', '52', '$332.28', '5330015847794', 'GASKET', 'Cameron', '9A.0199.0027C', 'M10', '2020 FEB 27'),
    @('9/20/2019', 'SPE5EJ-19-V-7230', '1', '$290.30', '5330010771949', 'GASKET', 'Munters', 'C60508-01', 'M31', '2020 FEB 27'),
    @('9/20/2019', 'SPE7M8-19-P-3507', '15', '$8,610.00', '5930011263228', 'SWITCH ASSEMBLY', 'Cole', '3600-3623', 'CP', '2020 MAR 16')
)

# First row that currently has no data (sheet currently ends at row 361).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($ws.Cells.Item(1,1).Value2 -eq $null) { $lastRow = 0 }
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 10))
    # Force Text format first so date-like / numeric-looking strings are
    # stored verbatim (matching the source data) instead of being parsed
    # into Excel dates/numbers.
    $rowRange.NumberFormat = "@"
    for ($c = 0; $c -lt 10; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $newRows[$i][$c]
    }
}

Write-Output ("Appended rows {0}-{1}" -f $startRow, ($startRow + $newRows.Count - 1))
